$d = $word.ActiveDocument

function New-ParaXmlFragment($bodyXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) The opening "Dear Hiring Manager," paragraph loses its greeting text
#    (becomes a blank paragraph), then three new paragraphs are inserted
#    right after it and before "I am excited to submit...": two blank ones
#    and a third one that reads "Dear Manager,".
# ---------------------------------------------------------------------------
$greetingPara = $d.Paragraphs(1)
$greetingRange = $greetingPara.Range
$textOnly = $d.Range($greetingRange.Start, $greetingRange.End - 1)
$textOnly.Delete()

$insertionPoint = $d.Range($greetingPara.Range.End - 1, $greetingPara.Range.End - 1)
$newParasBody = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p>' `
    + '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr></w:p>' `
    + '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr><w:r><w:t>Dear Manager,</w:t></w:r></w:p>'
$insertionPoint.InsertXML((New-ParaXmlFragment $newParasBody))

# ---------------------------------------------------------------------------
# 2) Merge the two runs "In regard to" + " the primary responsibilities..."
#    into a single run by replacing the whole sentence via Find/Replace.
# ---------------------------------------------------------------------------
$mergedText = "In regard to the primary responsibilities, I have experience handling confidential and sensitive materials in previous roles, and I am comfortable making deliveries and running errands on-campus. I am also excited about the opportunity to serve as a tour guide for Academic Hall and Academic Hall Dome, as I am knowledgeable about the history of Southeast Missouri State University."
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append four blank paragraphs (spacing 240) after "Chandrasekhar."
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$blank240 = ""
for ($i = 0; $i -lt 4; $i++) {
    $blank240 += '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr></w:p>'
}
$endPoint.InsertXML((New-ParaXmlFragment $blank240))
